$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(45, 8).Value = 201760
$ws.Cells.Item(45, 10).Value = 201760
$ws.Cells.Item(45, 12).Value = 605280
$ws.Cells.Item(45, 14).Value = -605664

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 4614.3335
$ws.Cells.Item(69, 9).Value = 3613
$ws.Cells.Item(69, 11).Value = 10839
$ws.Cells.Item(69, 13).Value = -9965

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(72, 8).Value = 4614.3335
$ws.Cells.Item(72, 9).Value = 3613
$ws.Cells.Item(72, 11).Value = 32517
$ws.Cells.Item(72, 13).Value = -28149

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 14766366
$ws.Cells.Item(112, 9).Value = 2650
$ws.Cells.Item(112, 10).Value = 22819302
$ws.Cells.Item(112, 11).Value = 7950
$ws.Cells.Item(112, 12).Value = 68457906
$ws.Cells.Item(112, 13).Value = -6842
$ws.Cells.Item(112, 14).Value = -68460122

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 4114.25
$ws.Cells.Item(116, 9).Value = 3637.1428
$ws.Cells.Item(116, 10).Value = 4371.154
$ws.Cells.Item(116, 11).Value = 3637.1428
$ws.Cells.Item(116, 12).Value = 4371.154
$ws.Cells.Item(116, 13).Value = -195.1428000000001
$ws.Cells.Item(116, 14).Value = -11255.154

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 5419.1113
$ws.Cells.Item(102, 9).Value = 2970
$ws.Cells.Item(102, 11).Value = 2970
$ws.Cells.Item(102, 13).Value = -1348

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 1316.7894
$ws.Cells.Item(110, 9).Value = 547.2593000000001
$ws.Cells.Item(110, 10).Value = 3205.6365
$ws.Cells.Item(110, 11).Value = 547.2593000000001
$ws.Cells.Item(110, 12).Value = 3205.6365
$ws.Cells.Item(110, 13).Value = 1497.7407
$ws.Cells.Item(110, 14).Value = -7295.636500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2059.3215
$ws.Cells.Item(20, 9).Value = 1645.8422
$ws.Cells.Item(20, 10).Value = 2932.2222
$ws.Cells.Item(20, 11).Value = 1645.8422
$ws.Cells.Item(20, 12).Value = 2932.2222
$ws.Cells.Item(20, 13).Value = -1398.8422
$ws.Cells.Item(20, 14).Value = -3426.2222

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 3589.2144
$ws.Cells.Item(99, 9).Value = 2924.9
$ws.Cells.Item(99, 11).Value = 2924.9
$ws.Cells.Item(99, 13).Value = -1426.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(132, 8).Value = 30000
$ws.Cells.Item(132, 10).Value = 30000
$ws.Cells.Item(132, 12).Value = 30000
$ws.Cells.Item(132, 14).Value = -40120

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2377.164
$ws.Cells.Item(134, 9).Value = 2201.7896
$ws.Cells.Item(134, 10).Value = 3376.8
$ws.Cells.Item(134, 11).Value = 6605.3688
$ws.Cells.Item(134, 12).Value = 10130.4
$ws.Cells.Item(134, 13).Value = -4070.3688
$ws.Cells.Item(134, 14).Value = -15200.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2884.19
$ws.Cells.Item(31, 9).Value = 1944.3
$ws.Cells.Item(31, 10).Value = 3510.7834
$ws.Cells.Item(31, 11).Value = 1944.3
$ws.Cells.Item(31, 12).Value = 3510.7834
$ws.Cells.Item(31, 13).Value = -1649.3
$ws.Cells.Item(31, 14).Value = -4100.7834

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2884.19
$ws.Cells.Item(34, 9).Value = 1944.3
$ws.Cells.Item(34, 10).Value = 3510.7834
$ws.Cells.Item(34, 11).Value = 1944.3
$ws.Cells.Item(34, 12).Value = 3510.7834
$ws.Cells.Item(34, 13).Value = -1742.3
$ws.Cells.Item(34, 14).Value = -3914.7834

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(74, 8).Value = 17342.908
$ws.Cells.Item(74, 10).Value = 17342.908
$ws.Cells.Item(74, 12).Value = 17342.908
$ws.Cells.Item(74, 14).Value = -19090.908

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(77, 8).Value = 17342.908
$ws.Cells.Item(77, 10).Value = 17342.908
$ws.Cells.Item(77, 12).Value = 52028.724
$ws.Cells.Item(77, 14).Value = -60764.724

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(100, 8).Value = 30000
$ws.Cells.Item(100, 10).Value = 30000
$ws.Cells.Item(100, 12).Value = 30000
$ws.Cells.Item(100, 14).Value = -32164

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 1073.2858
$ws.Cells.Item(107, 9).Value = 704.6667
$ws.Cells.Item(107, 10).Value = 1349.75
$ws.Cells.Item(107, 11).Value = 704.6667
$ws.Cells.Item(107, 12).Value = 1349.75
$ws.Cells.Item(107, 13).Value = 1215.3333
$ws.Cells.Item(107, 14).Value = -5189.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 34.23077
$ws.Cells.Item(2, 9).Value = 9.625
$ws.Cells.Item(2, 10).Value = 73.59999999999999
$ws.Cells.Item(2, 11).Value = 57.75
$ws.Cells.Item(2, 12).Value = 441.6
$ws.Cells.Item(2, 13).Value = 55.25
$ws.Cells.Item(2, 14).Value = -667.5999999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(64, 8).Value = 16284.571
$ws.Cells.Item(64, 9).Value = 1664
$ws.Cells.Item(64, 10).Value = 27250
$ws.Cells.Item(64, 11).Value = 4992
$ws.Cells.Item(64, 12).Value = 81750
$ws.Cells.Item(64, 13).Value = -4722
$ws.Cells.Item(64, 14).Value = -82290

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(67, 8).Value = 16284.571
$ws.Cells.Item(67, 9).Value = 1664
$ws.Cells.Item(67, 10).Value = 27250
$ws.Cells.Item(67, 11).Value = 4992
$ws.Cells.Item(67, 12).Value = 81750
$ws.Cells.Item(67, 13).Value = -4056
$ws.Cells.Item(67, 14).Value = -83622

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 839.01
$ws.Cells.Item(131, 9).Value = 392.1111
$ws.Cells.Item(131, 10).Value = 883.2088
$ws.Cells.Item(131, 11).Value = 1176.3333
$ws.Cells.Item(131, 12).Value = 2649.6264
$ws.Cells.Item(131, 13).Value = 3863.6667
$ws.Cells.Item(131, 14).Value = -12729.6264

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(141, 8).Value = 2510
$ws.Cells.Item(141, 10).Value = 6000
$ws.Cells.Item(141, 12).Value = 18000
$ws.Cells.Item(141, 14).Value = -28360

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4353.5293
$ws.Cells.Item(70, 9).Value = 4323.636
$ws.Cells.Item(70, 10).Value = 4408.3335
$ws.Cells.Item(70, 11).Value = 4323.636
$ws.Cells.Item(70, 12).Value = 4408.3335
$ws.Cells.Item(70, 13).Value = -4053.636
$ws.Cells.Item(70, 14).Value = -4948.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 4353.5293
$ws.Cells.Item(73, 9).Value = 4323.636
$ws.Cells.Item(73, 10).Value = 4408.3335
$ws.Cells.Item(73, 11).Value = 4323.636
$ws.Cells.Item(73, 12).Value = 4408.3335
$ws.Cells.Item(73, 13).Value = -3387.636
$ws.Cells.Item(73, 14).Value = -6280.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 56446.367
$ws.Cells.Item(102, 9).Value = 2882.6667
$ws.Cells.Item(102, 10).Value = 148269.86
$ws.Cells.Item(102, 11).Value = 2882.6667
$ws.Cells.Item(102, 12).Value = 148269.86
$ws.Cells.Item(102, 13).Value = -1260.6667
$ws.Cells.Item(102, 14).Value = -151513.86

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(106, 8).Value = 29166.666
$ws.Cells.Item(106, 10).Value = 29166.666
$ws.Cells.Item(106, 12).Value = 29166.666
$ws.Cells.Item(106, 14).Value = -31690.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(116, 8).Value = 29953.334
$ws.Cells.Item(116, 10).Value = 29953.334
$ws.Cells.Item(116, 12).Value = 29953.334
$ws.Cells.Item(116, 14).Value = -39131.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3779.6365
$ws.Cells.Item(132, 9).Value = 2861.3333
$ws.Cells.Item(132, 10).Value = 4415.385
$ws.Cells.Item(132, 11).Value = 8583.999899999999
$ws.Cells.Item(132, 12).Value = 13246.155
$ws.Cells.Item(132, 13).Value = -6053.999899999999
$ws.Cells.Item(132, 14).Value = -18306.155

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(13, 8).Value = 83337
$ws.Cells.Item(13, 9).Value = 49999
$ws.Cells.Item(13, 11).Value = 49999
$ws.Cells.Item(13, 13).Value = -49859

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(104, 8).Value = 34900
$ws.Cells.Item(104, 10).Value = 34900
$ws.Cells.Item(104, 12).Value = 34900
$ws.Cells.Item(104, 14).Value = -41888

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1524.8182
$ws.Cells.Item(113, 9).Value = 197
$ws.Cells.Item(113, 11).Value = 591
$ws.Cells.Item(113, 13).Value = 1579
